# Apply PnP table refresh: add LED1 and R1 rows, correct a handful of
# Mid-Y (and one Mid-X) values, rename the "Mayer" column to "Layer",
# resize the query table / defined name to the new A1:E30 extent, and
# move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a few existing Mid X / Mid Y values that changed on refresh ---
$ws.Range("C3").Value = 4            # C10
$ws.Range("B5").Value = 34           # C12
$ws.Range("C5").Value = 9            # C12
$ws.Range("C9").Value = 4            # C33
$ws.Range("C10").Value = 4           # C34
$ws.Range("C11").Value = 4           # C35
$ws.Range("C12").Value = 4           # C36

# --- Insert the new LED1 row right after L1 (row 19) ---
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "LED1"
$ws.Range("B20").Value = 36
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = "top"
$ws.Range("E20").Value = 180

# --- Insert the new R1 row right after Q10 (now row 22) ---
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "R1"
$ws.Range("B23").Value = 36
$ws.Range("C23").Value = 17
$ws.Range("D23").Value = "top"
$ws.Range("E23").Value = 0

# --- Rename the "Mayer" header to "Layer" ---
$ws.Range("D1").Value = "Layer"

# --- Resize table / autofilter to the new A1:E30 extent ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E30"))

# --- Resize the hidden ExternalData_1 defined name to match ---
$wb.Names.Item("Sheet1!ExternalData_1").RefersTo = "=Sheet1!`$A`$1:`$E`$30"

# --- Move the active selection, as Excel would leave it after refresh ---
$ws.Range("E20").Select()
